$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # SignIn
$ws2 = $wb.Worksheets.Item(2)   # doSignIn
$ws3 = $wb.Worksheets.Item(3)   # doSaveUpdates
$ws4 = $wb.Worksheets.Item(4)   # doSearch

# ---------------------------------------------------------------------------
# Drop every pre-existing hyperlink on the two sign-in sheets before the
# cells are rewritten, so no stale mailto: relationships survive.
# ---------------------------------------------------------------------------
$null = $ws1.Range("A1").Hyperlinks.Delete()
$null = $ws2.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Cell values -- written in the exact order the new unique strings need to
# land in xl/sharedStrings.xml (interleaved across sheets) so the resulting
# shared string table lines up with the target workbook.
# ---------------------------------------------------------------------------
$ws1.Range("A3").Value = "toubouachefazia@gmail.com"
$ws1.Range("B3").Value = "Fazia@96us"

$ws3.Range("B1").Value = "StreetAddress"
$ws3.Range("C1").Value = "city"
$ws3.Range("A2").Value = "223-334-5589"
$ws3.Range("C3").Value = "westfield"
$ws3.Range("A3").Value = "209-998-7765"
$ws3.Range("B3").Value = "1233 Main st"
$ws3.Range("B2").Value = "1177 Ritchie dr"
$ws3.Range("A4").Value = "209-4435-2314"
$ws3.Range("B4").Value = "1277 Main st"
$ws3.Range("A5").Value = "916-657-4576"

$ws1.Range("A2").Value = "fazousara96@yahoo.com"
$ws1.Range("B2").Value = "Abc1996@us"
$ws1.Range("A4").Value = "fazousaradouni@gmail.com"
$ws1.Range("B4").Value = "Fabc1996@us"

$ws4.Range("A1").Value = "SearchTerm"
$ws4.Range("A4").Value = "1235 IN, Carmel"
$ws4.Range("A3").Value = "1277 IN, Westfield"
$ws4.Range("A2").Value = "46074 IN, Westfield"
$ws4.Range("A5").Value = "46240 IN, Nora"
$ws4.Range("A6").Value = "46069 IN,Sheridan"

$ws3.Range("C2").Value = "Carmel"
$ws3.Range("C4").Value = "Noblesville"
$ws3.Range("C5").Value = "Fishers"

# Sheet2 re-uses strings already introduced above (toubouachefazia..., Fazia@96us,
# fazousaradouni..., Fabc1996@us, fazousara96..., Abc1996@us) so it does not
# introduce any new shared-string entries.
$ws2.Range("A2").Value = "toubouachefazia@gmail.com"
$ws2.Range("B2").Value = "Fazia@96us"
$ws2.Range("A3").Value = "fazousaradouni@gmail.com"
$ws2.Range("B3").Value = "Fabc1996@us"
$ws2.Range("A4").Value = "fazousara96@yahoo.com"
$ws2.Range("B4").Value = "Abc1996@us"

# ---------------------------------------------------------------------------
# Hyperlinks -- added in the exact order required so the generated r:id
# numbering (rId1, rId2, ...) matches the target relationships.
# ---------------------------------------------------------------------------
$null = $ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:toubouachefazia@gmail.com")
$null = $ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:Fazia@96us")
$null = $ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:fazousara96@yahoo.com")
$null = $ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:Abc1996@us")
$null = $ws1.Hyperlinks.Add($ws1.Range("A4"), "mailto:fazousaradouni@gmail.com")
$null = $ws1.Hyperlinks.Add($ws1.Range("B4"), "mailto:Fabc1996@us")
$ws1.Range("A2:B4").Style = "Hyperlink"

$null = $ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:toubouachefazia@gmail.com")
$null = $ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:Fazia@96us")
$null = $ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:fazousaradouni@gmail.com")
$null = $ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:Fabc1996@us")
$null = $ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:fazousara96@yahoo.com")
$null = $ws2.Hyperlinks.Add($ws2.Range("B4"), "mailto:Abc1996@us")
$ws2.Range("A2:B4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Selections -- applied last, non-active sheets first, so that sheet1
# ("SignIn") ends up as the selected/active tab, matching the target.
# ---------------------------------------------------------------------------
$null = $ws2.Range("A3").Select()
$null = $ws3.Range("A2").Select()
$null = $ws4.Range("A6").Select()
$null = $ws1.Range("B3").Select()
